# Update "想去人数" (interest/want-to-go count) figures in column F across
# all four sheets, per the "gh-pages output generated at 456a3b4" refresh.

$wb = $excel.ActiveWorkbook

# --- Sheet: 展览 (Exhibitions) ---
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F4").Value = 319    # was 318
$ws.Range("F5").Value = 5751   # was 5733
$ws.Range("F7").Value = 9777   # was 9740
$ws.Range("F10").Value = 3902  # was 3887
$ws.Range("F18").Value = 111   # was 110
$ws.Range("F20").Value = 625   # was 623
$ws.Range("F21").Value = 3924  # was 3916
$ws.Range("F22").Value = 139   # was 138
$ws.Range("F24").Value = 5382  # was 5370
$ws.Range("F26").Value = 2129  # was 2124
$ws.Range("F28").Value = 362   # was 360
$ws.Range("F29").Value = 8023  # was 7992
$ws.Range("F31").Value = 5     # was 3
$ws.Range("F32").Value = 2208  # was 2207
$ws.Range("F33").Value = 2217  # was 2213
$ws.Range("F34").Value = 1337  # was 1336
$ws.Range("F35").Value = 1318  # was 1313
$ws.Range("F36").Value = 22    # was 21
$ws.Range("F38").Value = 276   # was 273
$ws.Range("F39").Value = 253   # was 251
$ws.Range("F43").Value = 180   # was 179
$ws.Range("F44").Value = 1348  # was 1345
$ws.Range("F45").Value = 2122  # was 2113
$ws.Range("F46").Value = 137   # was 135
$ws.Range("F47").Value = 230   # was 229

# --- Sheet: 演出 (Performances) ---
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F4").Value = 149    # was 148
$ws.Range("F11").Value = 126   # was 125
$ws.Range("F20").Value = 18    # was 17

# --- Sheet: 本地生活 (Local life) ---
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 586    # was 583
$ws.Range("F3").Value = 768    # was 761

# --- Sheet: 全部类型 (All types) ---
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 768    # was 761
$ws.Range("F5").Value = 319    # was 318
$ws.Range("F6").Value = 5751   # was 5733
$ws.Range("F8").Value = 3903   # was 3887
$ws.Range("F15").Value = 111   # was 110
$ws.Range("F16").Value = 149   # was 148
$ws.Range("F18").Value = 625   # was 623
$ws.Range("F19").Value = 3924  # was 3916
$ws.Range("F21").Value = 139   # was 138
$ws.Range("F23").Value = 5382  # was 5370
$ws.Range("F25").Value = 2129  # was 2124
$ws.Range("F27").Value = 362   # was 360
$ws.Range("F28").Value = 8023  # was 7992
$ws.Range("F30").Value = 2208  # was 2207
$ws.Range("F31").Value = 2217  # was 2213
$ws.Range("F32").Value = 1337  # was 1336
$ws.Range("F33").Value = 1318  # was 1313
$ws.Range("F35").Value = 276   # was 273
$ws.Range("F36").Value = 253   # was 251
$ws.Range("F40").Value = 180   # was 179
$ws.Range("F42").Value = 1348  # was 1345
$ws.Range("F44").Value = 2122  # was 2113
$ws.Range("F45").Value = 137   # was 135
$ws.Range("F46").Value = 230   # was 229
$ws.Range("F48").Value = 18    # was 17
